{"js": "// 1. Fix the \"01.\" import line: {#import ...} -> {@import ...}\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet importPara = null;\nlet endProjectPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"01. {#import\") === 0) {\n    importPara = paragraphs.items[i];\n  }\n  if (t.indexOf(\"23. {end project}\") === 0) {\n    endProjectPara = paragraphs.items[i];\n  }\n}\n\nif (importPara) {\n  const fixedText = importPara.text.replace(\"{#import\", \"{@import\");\n  importPara.insertText(fixedText, \"Replace\");\n}\n\n// 2. Add a new paragraph \"24. {@xmlFormat}\" right after the \"23. {end project}\" paragraph,\n//    matching its Courier New / 18pt / en-US formatting.\nif (endProjectPara) {\n  endProjectPara.insertParagraph(\"24. {@xmlFormat}\", \"After\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Fix the \"01.\" import line: {#import ...} -> {@import ...}\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n  $p = $paras.Item($i)\n  $t = $p.Range.Text\n  if ($t.StartsWith(\"01. {#import\")) {\n    $find = $p.Range.Find\n    $find.Execute(\"{#import\", $false, $false, $false, $false, $false, $true, 1, $false, \"{@import\", 2)\n    break\n  }\n}\n\n# 2. Add a new paragraph \"24. {@xmlFormat}\" right after the \"23. {end project}\" paragraph,\n#    inheriting its Courier New / 18pt / en-US formatting.\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n  $p = $paras.Item($i)\n  $t = $p.Range.Text\n  if ($t.StartsWith(\"23. {end project}\")) {\n    $p.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($i + 1)\n    $newPara.Range.Text = \"24. {@xmlFormat}\"\n    break\n  }\n}\n"}
